$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.354.38"
Set-TextValue $ws.Range("E2") "  -1.39%  "
Set-TextValue $ws.Range("D3") "1.831.48"
Set-TextValue $ws.Range("E3") "  -1.14%  "
Set-TextValue $ws.Range("E4") "  -1.11%  "
Set-TextValue $ws.Range("D5") "314.77"
Set-TextValue $ws.Range("E5") "  -1.75%  "
Set-TextValue $ws.Range("D6") "1.005"
Set-TextValue $ws.Range("E6") "  -0.81%  "
Set-TextValue $ws.Range("D7") "0.4267"
Set-TextValue $ws.Range("E7") "  -1.22%  "
Set-TextValue $ws.Range("E8") "  -2.20%  "
Set-TextValue $ws.Range("D9") "0.07257"
Set-TextValue $ws.Range("E9") "  -1.94%  "
Set-TextValue $ws.Range("E10") "  -1.74%  "
Set-TextValue $ws.Range("D11") "21.10"
Set-TextValue $ws.Range("E11") "  -3.04%  "
Set-TextValue $ws.Range("D12") "1.824.27"
Set-TextValue $ws.Range("E12") "  -1.72%  "
Set-TextValue $ws.Range("D13") "6.716"
Set-TextValue $ws.Range("E13") "  -0.56%  "
Set-TextValue $ws.Range("D14") "0.07127"
Set-TextValue $ws.Range("E14") "  +0.04%  "
Set-TextValue $ws.Range("D15") "5.318"
Set-TextValue $ws.Range("E15") "  -3.12%  "
Set-TextValue $ws.Range("D16") "89.11"
Set-TextValue $ws.Range("E16") "  +0.71%  "
Set-TextValue $ws.Range("D17") "1.007"
Set-TextValue $ws.Range("E17") "  -0.97%  "
Set-TextValue $ws.Range("D18") "0.000008869"
Set-TextValue $ws.Range("E18") "  -1.80%  "
Set-TextValue $ws.Range("D19") "1.005"
Set-TextValue $ws.Range("E19") "  -0.85%  "
Set-TextValue $ws.Range("D20") "15.09"
Set-TextValue $ws.Range("E20") "  -2.70%  "
Set-TextValue $ws.Range("D21") "27.345.12"
Set-TextValue $ws.Range("E21") "  -1.51%  "
Set-TextValue $ws.Range("D22") "5.145"
Set-TextValue $ws.Range("E22") "  -2.40%  "
Set-TextValue $ws.Range("D23") "10.91"
Set-TextValue $ws.Range("E23") "  -2.32%  "
Set-TextValue $ws.Range("D24") "2.049.45"
Set-TextValue $ws.Range("E24") "  -2.12%  "
Set-TextValue $ws.Range("D25") "2.006"
Set-TextValue $ws.Range("E25") "  -1.19%  "
Set-TextValue $ws.Range("D26") "152.71"
Set-TextValue $ws.Range("E26") "  -2.26%  "
Set-TextValue $ws.Range("D27") "2.187"
Set-TextValue $ws.Range("E27") "  +6.39%  "
Set-TextValue $ws.Range("D28") "18.42"
Set-TextValue $ws.Range("E28") "  -1.35%  "
Set-TextValue $ws.Range("D29") "5.257"
Set-TextValue $ws.Range("E29") "  -3.16%  "
Set-TextValue $ws.Range("D30") "116.55"
Set-TextValue $ws.Range("E30") "  -4.34%  "
Set-TextValue $ws.Range("D31") "0.08900"
Set-TextValue $ws.Range("E31") "  -0.77%  "
Set-TextValue $ws.Range("D32") "1.204"
Set-TextValue $ws.Range("E32") "  -2.79%  "
Set-TextValue $ws.Range("D33") "0.7600"
Set-TextValue $ws.Range("E33") "  -2.36%  "
Set-TextValue $ws.Range("D34") "4.466"
Set-TextValue $ws.Range("E34") "  -2.42%  "
Set-TextValue $ws.Range("D35") "2.824"
Set-TextValue $ws.Range("E35") "  -3.52%  "
Set-TextValue $ws.Range("D36") "1.006"
Set-TextValue $ws.Range("E36") "  -0.80%  "
Set-TextValue $ws.Range("D37") "1.117"
Set-TextValue $ws.Range("E37") "  -2.48%  "
Set-TextValue $ws.Range("D38") "0.01986"
Set-TextValue $ws.Range("E38") "  +0.70%  "
Set-TextValue $ws.Range("D39") "0.05289"
Set-TextValue $ws.Range("E39") "  -0.70%  "
Set-TextValue $ws.Range("D40") "7.219"
Set-TextValue $ws.Range("E40") "  +2.81%  "
Set-TextValue $ws.Range("E41") "  +0.26%  "
Set-TextValue $ws.Range("E42") "  +1.05%  "
Set-TextValue $ws.Range("D43") "0.5086"
Set-TextValue $ws.Range("E43") "  -2.13%  "
Set-TextValue $ws.Range("D44") "8.718"
Set-TextValue $ws.Range("E44") "  -1.13%  "
Set-TextValue $ws.Range("D45") "10.64"
Set-TextValue $ws.Range("E45") "  -0.72%  "
Set-TextValue $ws.Range("D46") "107.94"
Set-TextValue $ws.Range("E46") "  -2.60%  "
Set-TextValue $ws.Range("D47") "0.4775"
Set-TextValue $ws.Range("E47") "  +0.69%  "
Set-TextValue $ws.Range("D48") "1.006"
Set-TextValue $ws.Range("E48") "  -0.78%  "
Set-TextValue $ws.Range("D49") "0.06388"
Set-TextValue $ws.Range("E49") "  -2.06%  "
Set-TextValue $ws.Range("D50") "1.670"
Set-TextValue $ws.Range("E50") "  -2.37%  "
Set-TextValue $ws.Range("D51") "1.855"
Set-TextValue $ws.Range("E51") "  -1.99%  "
